$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "NA" values under duplicate_image_filename (column E) for rows 2 through 21
for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 5).Value = "NA"
}
